$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 74 - this shifts the existing rows 74-77 down to 75-78
$ws.Rows.Item(74).Insert()

# Populate the new row 74 with this week's data point
$ws.Cells.Item(74, 1).Value = 11
$ws.Cells.Item(74, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(74, 3).Value = "Bíobío"
$ws.Cells.Item(74, 4).Value = 45021
$ws.Cells.Item(74, 5).Value = 8
$ws.Cells.Item(74, 6).Value = 100112031
$ws.Cells.Item(74, 7).Value = "Poroto verde"
$ws.Cells.Item(74, 8).Value = "Magnum"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 100
$ws.Cells.Item(74, 11).Value = 15000
$ws.Cells.Item(74, 12).Value = 16000
$ws.Cells.Item(74, 13).Value = 15500
$ws.Cells.Item(74, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(74, 15).Value = "Región Metropolitana"
$ws.Cells.Item(74, 16).Value = 620
$ws.Cells.Item(74, 17).Value = 25
$ws.Cells.Item(74, 18).Value = "Hortaliza"
